$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append rows 5-9 ---
$proximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "16:43:45", "16:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "16:43:45", "16:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "16:44:00", "16:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "16:44:26", "16:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "16:44:34", "16:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)

$startRow = 5
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $r = $startRow + $i
    $row = $proximityRows[$i]
    # Force column A to plain text so the literal date-like string
    # ("2026-02-01") is not auto-converted into a date serial number.
    $proximity.Range("A" + $r).NumberFormat = "@"
    $proximity.Cells.Item($r, 1).Value = $row[0]
    $proximity.Cells.Item($r, 2).Value = $row[1]
    $proximity.Cells.Item($r, 3).Value = $row[2]
    $proximity.Cells.Item($r, 4).Value = $row[3]
    $proximity.Cells.Item($r, 5).Value = $row[4]
    $proximity.Cells.Item($r, 6).Value = $row[5]
}

# --- Camera sheet: append rows 5-8 ---
$camera = $wb.Worksheets.Item("Camera")

$cameraRows = @(
    @("2026-02-01", "16:43:45", "16:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "16:44:01", "16:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "16:44:26", "16:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "16:44:35", "16:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow = 5
for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $r = $startRow + $i
    $row = $cameraRows[$i]
    # Force column A to plain text so the literal date-like string
    # ("2026-02-01") is not auto-converted into a date serial number.
    $camera.Range("A" + $r).NumberFormat = "@"
    $camera.Cells.Item($r, 1).Value = $row[0]
    $camera.Cells.Item($r, 2).Value = $row[1]
    $camera.Cells.Item($r, 3).Value = $row[2]
    $camera.Cells.Item($r, 4).Value = $row[3]
    $camera.Cells.Item($r, 5).Value = $row[4]
    $camera.Cells.Item($r, 6).Value = $row[5]
}
